# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the newly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (first worksheet) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 7005
$wsExhibit.Range("F3").Value = 55
$wsExhibit.Range("F5").Value = 87
$wsExhibit.Range("F6").Value = 1080
$wsExhibit.Range("F7").Value = 172
$wsExhibit.Range("F8").Value = 10

# --- Sheet "全部类型" (fourth worksheet) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 7005
$wsAll.Range("F3").Value = 55
$wsAll.Range("F5").Value = 87
$wsAll.Range("F6").Value = 1080
$wsAll.Range("F7").Value = 172
$wsAll.Range("F9").Value = 10
